$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column F is "dSF". Update the specific cells per the diff.
$ws.Range("F8").Value = 3
$ws.Range("F10").Value = -5
$ws.Range("F19").Value = -5
$ws.Range("F21").Value = 0
